
$d = $word.ActiveDocument

# Step 1: create a scratch paragraph at the end, apply default numbering to mint
# numId=2 / abstractNumId=1 with the decimal/lowerLetter/lowerRoman pattern, then
# remove the scratch paragraph's text but keep list defs in numbering.xml.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$rng0 = $lastPara.Range
$rng0.Collapse(0)
$rng0.InsertParagraphAfter()
$rng0.Collapse(0)
$scratchPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$scratchPara.Range.Text = "SCRATCH"
$scratchPara.Style = "ListParagraph"
$scratchPara.Range.ListFormat.ApplyNumberDefault()
$mintedNumId = $scratchPara.Range.ListFormat.List.ListID
Write-Host "minted numId: $mintedNumId"
$lt = $scratchPara.Range.ListFormat.ListTemplate
$lt.ListLevels.Item(2).NumberStyle = 4
$lt.ListLevels.Item(3).NumberStyle = 2
$lt.ListLevels.Item(5).NumberStyle = 4
$lt.ListLevels.Item(6).NumberStyle = 2
$lt.ListLevels.Item(8).NumberStyle = 4
$lt.ListLevels.Item(9).NumberStyle = 2

# Step 2: remove the scratch paragraph entirely (and its paragraph mark) so the
# document returns to its original last paragraph.
$delRng = $scratchPara.Range
$delRng.Delete()
Write-Host "Paragraphs after cleanup: $($d.Paragraphs.Count)"

# Step 3: insert the real body content, which already references numId=2 inline.
$rng = $d.Content
$rng.Collapse(0)
$xmlToInsert = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>10/10/2024 – The Beginnings of a Plan</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t>Install C++ compiler and build a simple hello world application to make sure the toolchain is squared away.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Define the data structures for the E-Vtol parent class, and one specific vehicle type child class. Let’s just start with the </w:t></w:r><w:r><w:t>Alpha Company.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Create</w:t></w:r><w:r><w:t xml:space="preserve"> one instance of an Alpha aircraft and perform some tests!</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Will need the main loop to step through time. </w:t></w:r><w:r><w:t xml:space="preserve">Can initially loop the corresponding </w:t></w:r><w:r><w:t>number</w:t></w:r><w:r><w:t xml:space="preserve"> of times for 1 simulation hour. This will be 60 times if the desired loop step is 1 minute.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Start in the idle state, move to the flying state, fly until the battery is dead, and then return to the idle state.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Record metrics during the flight, and verify that everything looks correct at the end!</w:t></w:r><w:r><w:t xml:space="preserve"> At this point in time, we’re not worrying about charging or faults. Data such as flight time, distance traveled, and total passenger miles should be easy to test now.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">The expected flight time for the Alpha aircraft is </w:t></w:r><w:r><w:t>100 minutes using the equation calculated yesterday. If the loop is kept at 1 hour for this first test, the Alpha plane should still be in the flying state at the end of the loop, with 60 minutes of recorded flight time. This will be an excellent test to perform!</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>If all passes, then loop for 2 hours, and expect the state to be idle at the end with 100 minutes of flight time.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Now we have an aircraft that can fly until its battery dies, and then stop flying. It would be a good time to incorporate one charger for our one aircraft to recharge with, and record its total time spent charging</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>With only one charger, we don’t have to worry about creating the queue right away.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>When the battery dies, the Alpha plane can be placed directly into the charging state.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>It will take 36 minutes to recharge its battery.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>For our 1-hour test, there should be no recharging.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">With the 2-hour test, the </w:t></w:r><w:r><w:t>total time spent charging should be 20 minutes, and the final state will be charging at the end of the loop.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>It’d be good to now perform a 3-hour test. What will we expect to see then?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>With one vehicle flying, recharging, and flying again; we’re at a very good starting point. Now will be a good time to incorporate another vehicle and another charger. How will we have to modify the code to implement this?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xmlToInsert)
Write-Host "Paragraphs: $($d.Paragraphs.Count)"
